# Apply updated Seasonality Index (column L) and Inventory Coverage (column H)
# values on the "Forecast Comparison" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Column L (Seasonality Index) updates
$ws.Range("L2").Value = 1.2
$ws.Range("L3").Value = 1.16
$ws.Range("L4").Value = 0.87
$ws.Range("L5").Value = 1
$ws.Range("L6").Value = 0.96
$ws.Range("L7").Value = 0.99
$ws.Range("L8").Value = 0.9
$ws.Range("L9").Value = 0.84
$ws.Range("L10").Value = 1.12
$ws.Range("L11").Value = 0.92
$ws.Range("L12").Value = 0.85
$ws.Range("L13").Value = 0.84
$ws.Range("L14").Value = 0.87
$ws.Range("L15").Value = 0.87
$ws.Range("L16").Value = 1.13
$ws.Range("L17").Value = 0.8100000000000001

# Column H (Inventory Coverage) updates
$ws.Range("H11").Value = 20.64
$ws.Range("H12").Value = 21.85
